{"js": "// Title Sheet \"READability\" edits:\n//   - \"LONG UH\"  -> \"Full Title\"\n//   - \"HILO\"     -> \"Address\"      (the standalone heading, not the \"UH HILO\" ones)\n//   - \"169\"      -> \"Section\"\n//   - \"09/10/2015\" -> \"09/11/2015\"\n//   - \"UH HILO\"  -> \"Short Title\"  (every occurrence; the last one also has a\n//                    trailing space run + tab run that collapse away)\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const text = para.text;\n\n  if (text === \"LONG UH\") {\n    para.getRange().insertText(\"Full Title\", \"Replace\");\n  } else if (text === \"HILO\") {\n    para.getRange().insertText(\"Address\", \"Replace\");\n  } else if (text === \"169\") {\n    para.getRange().insertText(\"Section\", \"Replace\");\n  } else if (text === \"09/10/2015\") {\n    para.getRange().insertText(\"09/11/2015\", \"Replace\");\n  } else if (text === \"UH HILO\" || text === \"UH HILO \\t\") {\n    // Replacing the whole-paragraph range collapses any extra trailing\n    // runs (the stray \" \" + tab runs on the final occurrence) into the\n    // single new text run, matching the other \"Short Title\" headings.\n    para.getRange().insertText(\"Short Title\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Title Sheet \"READability\" edits:\n#   - \"LONG UH\"  -> \"Full Title\"\n#   - \"HILO\"     -> \"Address\"      (the standalone heading, not the \"UH HILO\" ones)\n#   - \"169\"      -> \"Section\"\n#   - \"09/10/2015\" -> \"09/11/2015\"\n#   - \"UH HILO\"  -> \"Short Title\"  (every occurrence; the last one also has a\n#                    trailing space run + tab run that collapse away)\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $r = $p.Range\n    # Exclude the trailing paragraph mark so the comparison/replacement only\n    # ever touches the paragraph's visible content.\n    [void]$r.MoveEnd(1, -1)\n    $text = $r.Text\n\n    if ($text -eq \"LONG UH\") {\n        $r.Text = \"Full Title\"\n    } elseif ($text -eq \"HILO\") {\n        $r.Text = \"Address\"\n    } elseif ($text -eq \"169\") {\n        $r.Text = \"Section\"\n    } elseif ($text -eq \"09/10/2015\") {\n        $r.Text = \"09/11/2015\"\n    } elseif ($text -eq \"UH HILO\" -or $text -eq (\"UH HILO \" + [char]9)) {\n        # Setting .Text on the whole-paragraph range collapses any extra\n        # trailing runs (the stray \" \" + tab runs on the final occurrence)\n        # into the single new text run, matching the other headings.\n        $r.Text = \"Short Title\"\n    }\n}\n"}
